$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.991866
$ws.Range("H2").Value = 11.975598
$ws.Range("I2").Value = 0.2672133716845888
$ws.Range("J2").Value = 0.2672133716845888
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 73.19179533333333
$ws.Range("N2").Value = 219.575386
$ws.Range("O2").Value = 0.4454729128883617
$ws.Range("P2").Value = 0.4454729128883617
$ws.Range("Q2").Value = 292.171839270092
$ws.Range("R2").Value = 2629.546553430828
$ws.Range("S2").Value = 0.1190363190470542
$ws.Range("T2").Value = 0.1190363190470543
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.991866
$ws.Range("H3").Value = 11.975598
$ws.Range("I3").Value = 0.2672133716845888
$ws.Range("J3").Value = 0.2672133716845888
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.224257666666666
$ws.Range("N3").Value = 21.672773
$ws.Range("O3").Value = 0.04396956095378667
$ws.Range("P3").Value = 0.04396956095378668
$ws.Range("Q3").Value = 28.838268554806
$ws.Range("R3").Value = 259.544416993254
$ws.Range("S3").Value = 0.01174925463395238
$ws.Range("T3").Value = 0.01174925463395238
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.991866
$ws.Range("H4").Value = 11.975598
$ws.Range("I4").Value = 0.2672133716845888
$ws.Range("J4").Value = 0.2672133716845888
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 81.09049733333333
$ws.Range("N4").Value = 243.271492
$ws.Range("O4").Value = 0.4935473968103956
$ws.Range("P4").Value = 0.4935473968103956
$ws.Range("Q4").Value = 323.702399228024
$ws.Range("R4").Value = 2913.321593052216
$ws.Range("S4").Value = 0.1318824639878574
$ws.Range("T4").Value = 0.1318824639878575
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.991866
$ws.Range("H5").Value = 11.975598
$ws.Range("I5").Value = 0.2672133716845888
$ws.Range("J5").Value = 0.2672133716845888
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.794787
$ws.Range("N5").Value = 8.384361
$ws.Range("O5").Value = 0.01701012934745599
$ws.Range("P5").Value = 0.01701012934745599
$ws.Range("Q5").Value = 11.156415202542
$ws.Range("R5").Value = 100.407736822878
$ws.Range("S5").Value = 0.004545334015724689
$ws.Range("T5").Value = 0.00454533401572469
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.046736666666666
$ws.Range("H6").Value = 9.14021
$ws.Range("I6").Value = 0.2039469203963923
$ws.Range("J6").Value = 0.2039469203963923
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 73.19179533333333
$ws.Range("N6").Value = 219.575386
$ws.Range("O6").Value = 0.4454729128883617
$ws.Range("P6").Value = 0.4454729128883617
$ws.Range("Q6").Value = 222.9961265412289
$ws.Range("R6").Value = 2006.96513887106
$ws.Range("S6").Value = 0.0908528287035917
$ws.Range("T6").Value = 0.09085282870359171
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.046736666666666
$ws.Range("H7").Value = 9.14021
$ws.Range("I7").Value = 0.2039469203963923
$ws.Range("J7").Value = 0.2039469203963923
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.224257666666666
$ws.Range("N7").Value = 21.672773
$ws.Range("O7").Value = 0.04396956095378667
$ws.Range("P7").Value = 0.04396956095378668
$ws.Range("Q7").Value = 22.01041072248111
$ws.Range("R7").Value = 198.09369650233
$ws.Range("S7").Value = 0.008967456547706249
$ws.Range("T7").Value = 0.008967456547706252
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.046736666666666
$ws.Range("H8").Value = 9.14021
$ws.Range("I8").Value = 0.2039469203963923
$ws.Range("J8").Value = 0.2039469203963923
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 81.09049733333333
$ws.Range("N8").Value = 243.271492
$ws.Range("O8").Value = 0.4935473968103956
$ws.Range("P8").Value = 0.4935473968103956
$ws.Range("Q8").Value = 247.0613915437022
$ws.Range("R8").Value = 2223.55252389332
$ws.Range("S8").Value = 0.1006574716491364
$ws.Range("T8").Value = 0.1006574716491364
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.046736666666666
$ws.Range("H9").Value = 9.14021
$ws.Range("I9").Value = 0.2039469203963923
$ws.Range("J9").Value = 0.2039469203963923
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.794787
$ws.Range("N9").Value = 8.384361
$ws.Range("O9").Value = 0.01701012934745599
$ws.Range("P9").Value = 0.01701012934745599
$ws.Range("Q9").Value = 8.514980028423333
$ws.Range("R9").Value = 76.63482025581
$ws.Range("S9").Value = 0.003469163495957943
$ws.Range("T9").Value = 0.003469163495957944
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.337145333333333
$ws.Range("H10").Value = 10.011436
$ws.Range("I10").Value = 0.2233867209774804
$ws.Range("J10").Value = 0.2233867209774804
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 73.19179533333333
$ws.Range("N10").Value = 219.575386
$ws.Range("O10").Value = 0.4454729128883617
$ws.Range("P10").Value = 0.4454729128883617
$ws.Range("Q10").Value = 244.2516582349218
$ws.Range("R10").Value = 2198.264924114296
$ws.Range("S10").Value = 0.09951273329441788
$ws.Range("T10").Value = 0.0995127332944179
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.337145333333333
$ws.Range("H11").Value = 10.011436
$ws.Range("I11").Value = 0.2233867209774804
$ws.Range("J11").Value = 0.2233867209774804
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 7.224257666666666
$ws.Range("N11").Value = 21.672773
$ws.Range("O11").Value = 0.04396956095378667
$ws.Range("P11").Value = 0.04396956095378668
$ws.Range("Q11").Value = 24.10839775911422
$ws.Range("R11").Value = 216.975579832028
$ws.Range("S11").Value = 0.009822216044285861
$ws.Range("T11").Value = 0.009822216044285863
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.337145333333333
$ws.Range("H12").Value = 10.011436
$ws.Range("I12").Value = 0.2233867209774804
$ws.Range("J12").Value = 0.2233867209774804
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 81.09049733333333
$ws.Range("N12").Value = 243.271492
$ws.Range("O12").Value = 0.4935473968103956
$ws.Range("P12").Value = 0.4935473968103956
$ws.Range("Q12").Value = 270.6107747536124
$ws.Range("R12").Value = 2435.496972782512
$ws.Range("S12").Value = 0.1102519346204456
$ws.Range("T12").Value = 0.1102519346204456
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.337145333333333
$ws.Range("H13").Value = 10.011436
$ws.Range("I13").Value = 0.2233867209774804
$ws.Range("J13").Value = 0.2233867209774804
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.794787
$ws.Range("N13").Value = 8.384361
$ws.Range("O13").Value = 0.01701012934745599
$ws.Range("P13").Value = 0.01701012934745599
$ws.Range("Q13").Value = 9.326610394710666
$ws.Range("R13").Value = 83.939493552396
$ws.Range("S13").Value = 0.003799837018331002
$ws.Range("T13").Value = 0.003799837018331003
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.563122666666667
$ws.Range("H14").Value = 13.689368
$ws.Range("I14").Value = 0.3054529869415386
$ws.Range("J14").Value = 0.3054529869415386
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 73.19179533333333
$ws.Range("N14").Value = 219.575386
$ws.Range("O14").Value = 0.4454729128883617
$ws.Range("P14").Value = 0.4454729128883617
$ws.Range("Q14").Value = 333.9831402995609
$ws.Range("R14").Value = 3005.848262696048
$ws.Range("S14").Value = 0.1360710318432979
$ws.Range("T14").Value = 0.1360710318432979
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.563122666666667
$ws.Range("H15").Value = 13.689368
$ws.Range("I15").Value = 0.3054529869415386
$ws.Range("J15").Value = 0.3054529869415386
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 7.224257666666666
$ws.Range("N15").Value = 21.672773
$ws.Range("O15").Value = 0.04396956095378667
$ws.Range("P15").Value = 0.04396956095378668
$ws.Range("Q15").Value = 32.96517390860711
$ws.Range("R15").Value = 296.686565177464
$ws.Range("S15").Value = 0.01343063372784218
$ws.Range("T15").Value = 0.01343063372784219
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.563122666666667
$ws.Range("H16").Value = 13.689368
$ws.Range("I16").Value = 0.3054529869415386
$ws.Range("J16").Value = 0.3054529869415386
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 81.09049733333333
$ws.Range("N16").Value = 243.271492
$ws.Range("O16").Value = 0.4935473968103956
$ws.Range("P16").Value = 0.4935473968103956
$ws.Range("Q16").Value = 370.0258864330062
$ws.Range("R16").Value = 3330.232977897056
$ws.Range("S16").Value = 0.1507555265529561
$ws.Range("T16").Value = 0.1507555265529561
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.563122666666667
$ws.Range("H17").Value = 13.689368
$ws.Range("I17").Value = 0.3054529869415386
$ws.Range("J17").Value = 0.3054529869415386
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.794787
$ws.Range("N17").Value = 8.384361
$ws.Range("O17").Value = 0.01701012934745599
$ws.Range("P17").Value = 0.01701012934745599
$ws.Range("Q17").Value = 12.75295590820533
$ws.Range("R17").Value = 114.776603173848
$ws.Range("S17").Value = 0.005195794817442356
$ws.Range("T17").Value = 0.005195794817442357
